$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B100").Value = "8c38adc983273412ce5a788bae1cd763"
$ws.Range("B104").Value = "70ed1002822acd7a8d4a883eddad7c5a"
$ws.Range("B113").Value = "0d18b58d7bb14ac8516d60d305d3044a"
$ws.Range("B122").Value = "0a006fcf601c0dc3688178a86448c92d"
$ws.Range("B164").Value = "ce23e1e88a77288a83344f31e94882e0"
$ws.Range("B230").Value = "b606152ac708016201e949e62bd22efa"
$ws.Range("B233").Value = "4dc6992645510e489bbe6c13b9760931"
$ws.Range("B331").Value = "d9986ed4380897b50d61c0803314de7c"
$ws.Range("B342").Value = "987f8cbae45cd57dd33f7ec641011f88"
$ws.Range("B343").Value = "9c8e173b79f48d63f00af95644862e76"
$ws.Range("B419").Value = "afba4ee92bb44bede48ddf483ac24705"
$ws.Range("B619").Value = "2ce1ca607062c74c27dfcc4b74e09724"
$ws.Range("B623").Value = "1f0b3070b3e05c85fbf80ad9d3ccb14c"
$ws.Range("B628").Value = "b4c28e9a6e235253beea9f6a35999b21"
$ws.Range("B757").Value = "7f5feba9a45735fef77978b4a7635326"
$ws.Range("B760").Value = "9209ba79ce1d304c7323ecb2b6096f0b"
$ws.Range("B763").Value = "4d8faf6924deae6dcadf94b3c836e675"
$ws.Range("B767").Value = "bec68725ca3ed1d2d22a539f7a43ba56"
$ws.Range("B779").Value = "babf3fd530aff2ea45435a4292853ff1"
$ws.Range("B818").Value = "62404a0231c04bcfa99e99ab057a9cc0"
$ws.Range("B831").Value = "e6686e08a26163f0baac23e499746edf"
